# Apply the "commande.xlsx" update:
#  - row 11: "7000toggle" (interrupteur on/off) replaced by the new part
#    "D102J12S215PQA" with a price (4.03) and a hyperlink to its Digikey/CK page
#  - row 13: "TL1105T" (bouton) gets the same new hyperlink in column F
#  - row 26: quantity of "res 10k" increased from 6 to 10
#  - row 27: "res 2k" line item removed entirely (only the trailing "." in G stays)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPartUrl = "https://www.digikey.ca/en/products/detail/c-k/D102J12S215PQA/768267"

# --- Row 11: swap the component for the new "D102J12S215PQA" part ---
$ws.Range("A11").Value = "D102J12S215PQA"
$ws.Range("D11").Value = 4.03
$ws.Hyperlinks.Add($ws.Range("F11"), $newPartUrl, "", "", $newPartUrl)
$ws.Range("F11").Style = "Lien hypertexte"

# --- Row 13: add a hyperlink to the same new part's product page ---
$ws.Hyperlinks.Add($ws.Range("F13"), $newPartUrl, "", "", $newPartUrl)
$ws.Range("F13").Style = "Lien hypertexte"

# --- Row 26: bump the quantity of "res 10k" from 6 to 10 ---
$ws.Range("C26").Value = 10

# --- Row 27: remove the "res 2k" line item (keep only the G27 separator) ---
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("H27").ClearContents()

# --- Update the view: scroll up a bit and move the active selection ---
$ws.Range("H27").Select()

$wb.Application.Calculate()
